# Apply the commit's changes to the workbook:
#  - Add two new border-only cell styles (top+bottom, and top+bottom+right)
#    and apply them to the "spacer" cells of the merged group headers
#    (the cells between the first and last cell of each B1:D1 / E1:G1 merge).
#  - Rename the "fedcore" column header to "approach".
#  - Clear the stray empty cell G5 on the computational_comparison sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# xlPasteFormats
$xlPasteFormats = -4122
# xlLineStyleNone / continuous thin line
$xlLineStyleNone = -4142
$xlContinuous    = 1

# --- Build the two new border styles from scratch (on C1 / D1 of the first sheet) ---

# C1: reset to the plain/default format, then add a thin top+bottom border only.
$ws1.Range("B4").Copy()
$ws1.Range("C1").PasteSpecial($xlPasteFormats)
$midCell = $ws1.Range("C1")
$midCell.Borders.LineStyle = $xlContinuous           # all 4 edges thin
$midCell.Borders.Item(7).LineStyle = $xlLineStyleNone  # xlEdgeLeft
$midCell.Borders.Item(10).LineStyle = $xlLineStyleNone # xlEdgeRight

# D1: reset to the plain/default format, then add a thin top+bottom+right border.
$ws1.Range("B4").Copy()
$ws1.Range("D1").PasteSpecial($xlPasteFormats)
$endCell = $ws1.Range("D1")
$endCell.Borders.LineStyle = $xlContinuous           # all 4 edges thin
$endCell.Borders.Item(7).LineStyle = $xlLineStyleNone  # xlEdgeLeft

# --- Reuse those two freshly-minted styles on every other "spacer" cell ---
# (copying the already-resolved format avoids Excel from registering
#  intermediate/duplicate style records)

$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Rename "fedcore" header to "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty G5 cell on computational_comparison ---
$ws2.Range("G5").ClearContents()
